$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.47585300000000003
$ws.Range("C3").Value = 0.48208499999999999
$ws.Range("C4").Value = 0.47908600000000001
$ws.Range("C5").Value = 0.47655199999999998

$ws.Range("I5").Select()
